$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 21.6679573059082
$ws.Range("C3").Value = 20.98608016967773
$ws.Range("C4").Value = 20.47920227050781
$ws.Range("C5").Value = 20.65205574035645
$ws.Range("C6").Value = 23.3771800994873
